# Apply the edit described by the diff:
#  - replace the rich-text shared string used in column A (rows 2-11)
#    with the plain string "D0(mu)"
#  - change the sheet's active/selected cell from I14 to C15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "D0(mu)"
}

$ws.Range("C15").Select()
